$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replicate the formatting of the last existing data row (row 5) down to the
# four new rows (6-9) so the new rows pick up the same column styles
# (date format in A, text format in B, centered text in C).
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E9").PasteSpecial(-4122)

# Row 6: 2014-03-12, Change 04, JEB, System Design document updated, Done
$ws.Range("A6").Value = 41710
$ws.Range("B6").Value = "04"
$ws.Range("C6").Value = "JEB"
$ws.Range("D6").Value = "System Design document updated"
$ws.Range("E6").Value = "Done"

# Row 7: 2014-03-12, Change 05, JEB, Test Cases Document for Scheduler Updated, Done
$ws.Range("A7").Value = 41710
$ws.Range("B7").Value = "05"
$ws.Range("C7").Value = "JEB"
$ws.Range("D7").Value = "Test Cases Document for Scheduler Updated"
$ws.Range("E7").Value = "Done"

# Row 8: 2014-03-12, Change 06, JEB, Traceability document for Scheduler Updated, Done
$ws.Range("A8").Value = 41710
$ws.Range("B8").Value = "06"
$ws.Range("C8").Value = "JEB"
$ws.Range("D8").Value = "Traceability document for Scheduler Updated"
$ws.Range("E8").Value = "Done"

# Row 9: 2014-03-13, Change 07, JEB, Report Documents updated, Done
$ws.Range("A9").Value = 41711
$ws.Range("B9").Value = "07"
$ws.Range("C9").Value = "JEB"
$ws.Range("D9").Value = "Report Documents updated"
$ws.Range("E9").Value = "Done"

# Update selection to match the post-edit cursor position (next empty cell)
$ws.Range("E10").Select()
